# "removed some test cases from the run"
# Flip the Execute flag (column A) from "Y" to "N" for the five
# VisitDetailsModalTest scenarios (rows 4-8) on the Scenarios sheet, and
# leave the workbook with the Scenarios sheet active/selected at A2
# (instead of Parameters being the active tab).

$wb = $excel.ActiveWorkbook

$scenarios = $wb.Worksheets.Item("Scenarios")

# Disable ("N") the test cases in rows 4-8 (column A = Execute flag).
$scenarios.Range("A4:A8").Value = "N"

# Make "Scenarios" the active sheet/tab again, with A2 selected, which
# also clears the previous active-tab state on "Parameters".
$scenarios.Activate()
$scenarios.Range("A2").Select()
